# Apply updates to Inscritos/Pagos/Inscricoes homologadas columns (E, F, H)
# per the diff: 2025-1_GestaoResultado_ResumoInscricoes_Integrado.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E4").Value = 42
$ws.Range("E10").Value = 477
$ws.Range("F10").Value = 239
$ws.Range("H10").Value = 239
$ws.Range("E11").Value = 323
$ws.Range("F11").Value = 182
$ws.Range("H11").Value = 182
$ws.Range("E12").Value = 468
$ws.Range("F12").Value = 254
$ws.Range("H12").Value = 254
$ws.Range("E13").Value = 122
$ws.Range("F13").Value = 64
$ws.Range("H13").Value = 64
$ws.Range("E14").Value = 120
$ws.Range("E15").Value = 155
$ws.Range("E17").Value = 90
$ws.Range("F17").Value = 47
$ws.Range("H17").Value = 47
$ws.Range("E23").Value = 191
$ws.Range("E24").Value = 201
$ws.Range("F24").Value = 109
$ws.Range("H24").Value = 109
$ws.Range("E25").Value = 250
$ws.Range("F25").Value = 122
$ws.Range("H25").Value = 122
$ws.Range("E26").Value = 150
$ws.Range("F26").Value = 91
$ws.Range("H26").Value = 91
$ws.Range("E27").Value = 310
$ws.Range("E28").Value = 190
$ws.Range("E30").Value = 196
$ws.Range("E32").Value = 177
$ws.Range("F32").Value = 104
$ws.Range("H32").Value = 104
$ws.Range("E33").Value = 273
$ws.Range("E34").Value = 207
$ws.Range("F34").Value = 133
$ws.Range("H34").Value = 133
$ws.Range("E35").Value = 142
$ws.Range("E36").Value = 66
$ws.Range("F37").Value = 73
$ws.Range("H37").Value = 73
$ws.Range("E38").Value = 87
$ws.Range("E39").Value = 175
$ws.Range("E40").Value = 249
$ws.Range("F40").Value = 114
$ws.Range("H40").Value = 114
$ws.Range("F41").Value = 175
$ws.Range("H41").Value = 175
$ws.Range("E42").Value = 352
$ws.Range("F42").Value = 190
$ws.Range("H42").Value = 190
$ws.Range("E43").Value = 109
$ws.Range("E44").Value = 299
$ws.Range("F44").Value = 148
$ws.Range("H44").Value = 148
$ws.Range("E45").Value = 134
$ws.Range("E46").Value = 302
$ws.Range("E47").Value = 427
$ws.Range("E49").Value = 276
$ws.Range("E50").Value = 234
$ws.Range("F50").Value = 108
$ws.Range("H50").Value = 108
